$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CONDUCTOR_files")

# The dummy/placeholder external-file names are no longer used; replace them
# with the literal "none" (test data update for CASE_1_ITER_like_LTS).
$ws.Range("E8").Value  = "none"   # EXTERNAL_ALPHAB (was alphab_dummy.xlsx)
$ws.Range("E9").Value  = "none"   # EXTERNAL_BFIELD (was bfield.xlsx)
$ws.Range("E10").Value = "none"   # EXTERNAL_CURRENT (was I_file_dummy.xlsx)
$ws.Range("E11").Value = "none"   # EXTERNAL_FLOW (was flow_dummy.xlsx)
$ws.Range("E12").Value = "none"   # EXTERNAL_HEAT (was Q_file_dummy.xlsx)
$ws.Range("E13").Value = "none"   # EXTERNAL_STRAIN (was strain_dummy.xlsx)
$ws.Range("E15").Value = "none"   # EXTERNAL_GRID (was spatial_discretization.xlsx)

# Make CONDUCTOR_files the active/selected sheet with E15 selected, matching
# the saved view state in the target workbook.
$ws.Activate()
$ws.Range("E15").Select()
